$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.789.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +15.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.732.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9964"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9895"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.81%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.64%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +23.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.225"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07564"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9924"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.381"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.045"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.734.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001147"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9884"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06782"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "85.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +13.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.380"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "25.728.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.439"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.835"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "134.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.923.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.173"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +28.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.754"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.786"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08593"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06700"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.575"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.13%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02465"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.90%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.172"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.34%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.290"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2195"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6442"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +12.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9891"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6241"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.878"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.135"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07476"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.78%  "
